$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "nom" column (CL), shifting "nom" -> CM
# and "url_produit" -> CN. This mirrors a new price-check timestamp column
# being appended to the history, just before the descriptive columns.
$ws.Columns("CL").Insert()

# New column header: timestamp of this run
$ws.Range("CL1").Value = "2026-01-31 18:16:35"

# Populate the new column with the latest known price for each product,
# i.e. duplicate the previous last column (CK, now shifted before CL).
# Column CK = 89, new column CL = 90 (1-based).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ckValue = $ws.Cells.Item($r, 89).Value()
    if ($ckValue -ne "") {
        $ws.Cells.Item($r, 90).Value = $ckValue
    }
}
